$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3000.5
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H62").Value = 5174.5
$ws.Range("I62").Value = 2377
$ws.Range("K62").Value = 2377
$ws.Range("M62").Value = -1753
$ws.Range("H65").Value = 5174.5
$ws.Range("I65").Value = 2377
$ws.Range("K65").Value = 11885
$ws.Range("M65").Value = -8765
$ws.Range("H69").Value = 5666.6665
$ws.Range("J69").Value = 5666.6665
$ws.Range("L69").Value = 16999.9995
$ws.Range("N69").Value = -18747.9995
$ws.Range("H72").Value = 5666.6665
$ws.Range("J72").Value = 5666.6665
$ws.Range("L72").Value = 50999.9985
$ws.Range("N72").Value = -59735.9985
$ws.Range("H76").Value = 4742
$ws.Range("I76").Value = 4656
$ws.Range("K76").Value = 4656
$ws.Range("M76").Value = -4341
$ws.Range("H79").Value = 4742
$ws.Range("I79").Value = 4656
$ws.Range("K79").Value = 4656
$ws.Range("M79").Value = -3564
$ws.Range("H88").Value = 2942.375
$ws.Range("J88").Value = 3248.4285
$ws.Range("L88").Value = 3248.4285
$ws.Range("N88").Value = -4060.4285
$ws.Range("H91").Value = 2942.375
$ws.Range("J91").Value = 3248.4285
$ws.Range("L91").Value = 3248.4285
$ws.Range("N91").Value = -6056.4285
$ws.Range("H112").Value = 3155.6667
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3155.6667
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 9467.000100000001
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -11683.0001
$ws.Range("H132").Value = 6852.1113
$ws.Range("J132").Value = 6667
$ws.Range("L132").Value = 20001
$ws.Range("N132").Value = -25061
$ws.Range("H137").Value = 1612.125
$ws.Range("I137").Value = 1612.125
$ws.Range("K137").Value = 4836.375
$ws.Range("M137").Value = -2286.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 100925.29
$ws.Range("I21").Value = 140296.4
$ws.Range("K21").Value = 140296.4
$ws.Range("M21").Value = -139922.4
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 2484.5
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 2484.5
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -2784.5
$ws.Range("H36").Value = 10253.25
$ws.Range("I36").Value = 12254.333
$ws.Range("J36").Value = 4250
$ws.Range("K36").Value = 12254.333
$ws.Range("L36").Value = 4250
$ws.Range("M36").Value = -11908.333
$ws.Range("N36").Value = -4942
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 53590.332
$ws.Range("J74").Value = 53590.332
$ws.Range("L74").Value = 53590.332
$ws.Range("N74").Value = -55462.332
$ws.Range("H77").Value = 53590.332
$ws.Range("J77").Value = 53590.332
$ws.Range("L77").Value = 160770.996
$ws.Range("N77").Value = -170130.996
$ws.Range("H105").Value = 3106.3215
$ws.Range("I105").Value = 2879.08
$ws.Range("K105").Value = 2879.08
$ws.Range("M105").Value = -1132.08
$ws.Range("H107").Value = 555.6
$ws.Range("I107").Value = 501.5
$ws.Range("J107").Value = 772
$ws.Range("K107").Value = 501.5
$ws.Range("L107").Value = 772
$ws.Range("M107").Value = 1418.5
$ws.Range("N107").Value = -4612
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 935.9
$ws.Range("I16").Value = 928.7778
$ws.Range("K16").Value = 928.7778
$ws.Range("M16").Value = -641.7778
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 2439.5715
$ws.Range("I31").Value = 1763.75
$ws.Range("K31").Value = 1763.75
$ws.Range("M31").Value = -1468.75
$ws.Range("H34").Value = 2439.5715
$ws.Range("I34").Value = 1763.75
$ws.Range("K34").Value = 1763.75
$ws.Range("M34").Value = -1561.75
$ws.Range("H63").Value = 93329.664
$ws.Range("I63").Value = 80000
$ws.Range("K63").Value = 80000
$ws.Range("M63").Value = -79314
$ws.Range("H66").Value = 93329.664
$ws.Range("I66").Value = 80000
$ws.Range("K66").Value = 240000
$ws.Range("M66").Value = -236568
$ws.Range("H105").Value = 1431.6923
$ws.Range("I105").Value = 794.8570999999999
$ws.Range("J105").Value = 2174.6667
$ws.Range("K105").Value = 794.8570999999999
$ws.Range("L105").Value = 2174.6667
$ws.Range("M105").Value = 952.1429000000001
$ws.Range("N105").Value = -5668.6667
$ws.Range("H113").Value = 935.9
$ws.Range("I113").Value = 928.7778
$ws.Range("K113").Value = 928.7778
$ws.Range("M113").Value = 1241.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 764.6667
$ws.Range("J34").Value = 829
$ws.Range("L34").Value = 2487
$ws.Range("N34").Value = -2655
$ws.Range("H39").Value = 987.5
$ws.Range("J39").Value = 700
$ws.Range("L39").Value = 2100
$ws.Range("N39").Value = -2688
$ws.Range("H55").Value = 911
$ws.Range("J55").Value = 1125
$ws.Range("L55").Value = 3375
$ws.Range("N55").Value = -3729
$ws.Range("H103").Value = 3180
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 3180
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 9540
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -11298
$ws.Range("H114").Value = 903.125
$ws.Range("I114").Value = 1075
$ws.Range("K114").Value = 3225
$ws.Range("M114").Value = 29
$ws.Range("H132").Value = 1074.4
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1093
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 9837
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -14897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1658
$ws.Range("J80").Value = 1475
$ws.Range("L80").Value = 1475
$ws.Range("N80").Value = -3471
$ws.Range("H83").Value = 1658
$ws.Range("J83").Value = 1475
$ws.Range("L83").Value = 7375
$ws.Range("N83").Value = -17359
$ws.Range("H122").Value = 1618.875
$ws.Range("I122").Value = 1535.1538
$ws.Range("J122").Value = 1981.6666
$ws.Range("K122").Value = 4605.4614
$ws.Range("L122").Value = 5944.9998
$ws.Range("M122").Value = -2155.4614
$ws.Range("N122").Value = -10844.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1477.5
$ws.Range("J22").Value = 1553.5
$ws.Range("L22").Value = 1553.5
$ws.Range("N22").Value = -2143.5
$ws.Range("H27").Value = 1477.5
$ws.Range("J27").Value = 1553.5
$ws.Range("L27").Value = 1553.5
$ws.Range("N27").Value = -1767.5
$ws.Range("H46").Value = 2154.2104
$ws.Range("I46").Value = 1590
$ws.Range("K46").Value = 1590
$ws.Range("M46").Value = -1402
$ws.Range("H55").Value = 718.5625
$ws.Range("J55").Value = 799.0714
$ws.Range("L55").Value = 799.0714
$ws.Range("N55").Value = -1145.0714
$ws.Range("H82").Value = 3581.7273
$ws.Range("I82").Value = 1700
$ws.Range("J82").Value = 3769.9
$ws.Range("K82").Value = 1700
$ws.Range("L82").Value = 3769.9
$ws.Range("M82").Value = -1339
$ws.Range("N82").Value = -4491.9
$ws.Range("H85").Value = 3581.7273
$ws.Range("I85").Value = 1700
$ws.Range("J85").Value = 3769.9
$ws.Range("K85").Value = 1700
$ws.Range("L85").Value = 3769.9
$ws.Range("M85").Value = -452
$ws.Range("N85").Value = -6265.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1413.3334
$ws.Range("I8").Value = 750
$ws.Range("K8").Value = 750
$ws.Range("M8").Value = -610
$ws.Range("H94").Value = 23109.334
$ws.Range("J94").Value = 23109.334
$ws.Range("L94").Value = 23109.334
$ws.Range("N94").Value = -24911.334
$ws.Range("H97").Value = 15190.667
$ws.Range("J97").Value = 15190.667
$ws.Range("L97").Value = 15190.667
$ws.Range("N97").Value = -17172.667
